$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4577.6665
$ws.Range("I40").Value = 4144.6665
$ws.Range("J40").Value = 4794.1665
$ws.Range("K40").Value = 4144.6665
$ws.Range("L40").Value = 4794.1665
$ws.Range("M40").Value = -3969.6665
$ws.Range("N40").Value = -5144.1665
$ws.Range("H43").Value = 8041.778
$ws.Range("J43").Value = 8422
$ws.Range("L43").Value = 8422
$ws.Range("N43").Value = -8560
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H74").Value = 6815024.5
$ws.Range("I74").Value = 28576964
$ws.Range("K74").Value = 28576964
$ws.Range("M74").Value = -28576028
$ws.Range("H76").Value = 83340270
$ws.Range("I76").Value = 142863380
$ws.Range("J76").Value = 7932.4
$ws.Range("K76").Value = 142863380
$ws.Range("L76").Value = 7932.4
$ws.Range("M76").Value = -142863065
$ws.Range("N76").Value = -8562.4
$ws.Range("H77").Value = 6815024.5
$ws.Range("I77").Value = 28576964
$ws.Range("K77").Value = 142884820
$ws.Range("M77").Value = -142880140
$ws.Range("H79").Value = 83340270
$ws.Range("I79").Value = 142863380
$ws.Range("J79").Value = 7932.4
$ws.Range("K79").Value = 142863380
$ws.Range("L79").Value = 7932.4
$ws.Range("M79").Value = -142862288
$ws.Range("N79").Value = -10116.4
$ws.Range("H103").Value = 1630.75
$ws.Range("J103").Value = 2412.75
$ws.Range("L103").Value = 7238.25
$ws.Range("N103").Value = -8410.25

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 684.25
$ws.Range("I5").Value = 684.25
$ws.Range("K5").Value = 684.25
$ws.Range("M5").Value = -572.25
$ws.Range("H132").Value = 2605486.5
$ws.Range("I132").Value = 5051813.5
$ws.Range("K132").Value = 15155440.5
$ws.Range("M132").Value = -15152910.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 684.25
$ws.Range("I4").Value = 684.25
$ws.Range("K4").Value = 684.25
$ws.Range("M4").Value = -569.25
$ws.Range("H82").Value = 50000
$ws.Range("I82").Value = 50000
$ws.Range("K82").Value = 50000
$ws.Range("M82").Value = -49617
$ws.Range("H85").Value = 50000
$ws.Range("I85").Value = 50000
$ws.Range("K85").Value = 50000
$ws.Range("M85").Value = -48674
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 2750
$ws.Range("I23").Value = 2750
$ws.Range("K23").Value = 2750
$ws.Range("M23").Value = -2510
$ws.Range("H27").Value = 2750
$ws.Range("I27").Value = 2750
$ws.Range("K27").Value = 2750
$ws.Range("M27").Value = -2558
$ws.Range("H58").Value = 34487940
$ws.Range("I58").Value = 47622012
$ws.Range("K58").Value = 47622012
$ws.Range("M58").Value = -47621809
$ws.Range("H86").Value = 14105.65
$ws.Range("I86").Value = 14037.571
$ws.Range("J86").Value = 14264.5
$ws.Range("K86").Value = 14037.571
$ws.Range("L86").Value = 14264.5
$ws.Range("M86").Value = -12914.571
$ws.Range("N86").Value = -16510.5
$ws.Range("H89").Value = 14105.65
$ws.Range("I89").Value = 14037.571
$ws.Range("J89").Value = 14264.5
$ws.Range("K89").Value = 70187.855
$ws.Range("L89").Value = 71322.5
$ws.Range("M89").Value = -64571.855
$ws.Range("N89").Value = -82554.5
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 66674080
$ws.Range("I105").Value = 111116970
$ws.Range("K105").Value = 111116970
$ws.Range("M105").Value = -111115223
$ws.Range("H106").Value = 22500
$ws.Range("J106").Value = 22500
$ws.Range("L106").Value = 22500
$ws.Range("N106").Value = -25024
$ws.Range("H132").Value = 6558.52
$ws.Range("I132").Value = 5445.684
$ws.Range("K132").Value = 16337.052
$ws.Range("M132").Value = -13807.052
$ws.Range("H136").Value = 34487940
$ws.Range("I136").Value = 47622012
$ws.Range("K136").Value = 142866036
$ws.Range("M136").Value = -142863486
$ws.Range("H141").Value = 207798.2
$ws.Range("I141").Value = 415000
$ws.Range("J141").Value = 69663.664
$ws.Range("K141").Value = 415000
$ws.Range("L141").Value = 69663.664
$ws.Range("M141").Value = -409820
$ws.Range("N141").Value = -80023.664

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 627703.1
$ws.Range("I5").Value = 10000
$ws.Range("K5").Value = 30000
$ws.Range("M5").Value = -29888
$ws.Range("H74").Value = 20798.6
$ws.Range("J74").Value = 26655.666
$ws.Range("L74").Value = 79966.99800000001
$ws.Range("N74").Value = -82088.99800000001
$ws.Range("H77").Value = 20798.6
$ws.Range("J77").Value = 26655.666
$ws.Range("L77").Value = 239900.994
$ws.Range("N77").Value = -250508.994
$ws.Range("H81").Value = 2650.875
$ws.Range("I81").Value = 2188.3333
$ws.Range("J81").Value = 2928.4
$ws.Range("K81").Value = 6564.999899999999
$ws.Range("L81").Value = 8785.200000000001
$ws.Range("M81").Value = -5441.999899999999
$ws.Range("N81").Value = -11031.2
$ws.Range("H82").Value = 21749.75
$ws.Range("J82").Value = 28333
$ws.Range("L82").Value = 84999
$ws.Range("N82").Value = -85811
$ws.Range("H84").Value = 2650.875
$ws.Range("I84").Value = 2188.3333
$ws.Range("J84").Value = 2928.4
$ws.Range("K84").Value = 19694.9997
$ws.Range("L84").Value = 26355.6
$ws.Range("M84").Value = -14078.9997
$ws.Range("N84").Value = -37587.60000000001
$ws.Range("H85").Value = 21749.75
$ws.Range("J85").Value = 28333
$ws.Range("L85").Value = 84999
$ws.Range("N85").Value = -87807
$ws.Range("H87").Value = 19505.4
$ws.Range("I87").Value = 13722
$ws.Range("J87").Value = 33000
$ws.Range("K87").Value = 41166
$ws.Range("L87").Value = 99000
$ws.Range("M87").Value = -39918
$ws.Range("N87").Value = -101496
$ws.Range("H90").Value = 19505.4
$ws.Range("I90").Value = 13722
$ws.Range("J90").Value = 33000
$ws.Range("K90").Value = 123498
$ws.Range("L90").Value = 297000
$ws.Range("M90").Value = -117258
$ws.Range("N90").Value = -309480
$ws.Range("H135").Value = 627703.1
$ws.Range("I135").Value = 10000
$ws.Range("K135").Value = 90000
$ws.Range("M135").Value = -87465

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 9998
$ws.Range("J40").Value = 9998
$ws.Range("L40").Value = 9998
$ws.Range("N40").Value = -10300
$ws.Range("H47").Value = 15031
$ws.Range("J47").Value = 15031
$ws.Range("L47").Value = 15031
$ws.Range("N47").Value = -16167
$ws.Range("H55").Value = 20258
$ws.Range("I55").Value = 17999.5
$ws.Range("J55").Value = 22516.5
$ws.Range("K55").Value = 17999.5
$ws.Range("L55").Value = 22516.5
$ws.Range("M55").Value = -17672.5
$ws.Range("N55").Value = -23170.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2902.45
$ws.Range("I16").Value = 1144.875
$ws.Range("J16").Value = 9932.75
$ws.Range("K16").Value = 1144.875
$ws.Range("L16").Value = 9932.75
$ws.Range("M16").Value = -974.875
$ws.Range("N16").Value = -10272.75
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H40").Value = 4187
$ws.Range("I40").Value = 3655.7144
$ws.Range("K40").Value = 3655.7144
$ws.Range("M40").Value = -3519.7144
$ws.Range("H93").Value = 1123
$ws.Range("I93").Value = 1734.375
$ws.Range("J93").Value = 579.55554
$ws.Range("K93").Value = 1734.375
$ws.Range("L93").Value = 579.55554
$ws.Range("M93").Value = -486.375
$ws.Range("N93").Value = -3075.55554
$ws.Range("H122").Value = 5235.6665
$ws.Range("J122").Value = 6234.6
$ws.Range("L122").Value = 18703.8
$ws.Range("N122").Value = -23603.8

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 819.5
$ws.Range("I107").Value = 819.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2458.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -538.5
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 5335.0386
$ws.Range("I132").Value = 1920.875
$ws.Range("J132").Value = 10797.7
$ws.Range("K132").Value = 5762.625
$ws.Range("L132").Value = 32393.1
$ws.Range("M132").Value = -3232.625
$ws.Range("N132").Value = -37453.10000000001
